$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update contact details (second entry's first name, hobbies trailing commas
# cleaned up, and third/fourth edits to the second row of data)
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# Refresh the font color on the Pincode/Phone number cells to solid black
$ws.Range("I2").Font.Color = 0
$ws.Range("K2").Font.Color = 0
$ws.Range("I3").Font.Color = 0
$ws.Range("K3").Font.Color = 0

# Row heights grew slightly after the edit/re-save
$ws.Rows.Item(1).RowHeight = 19.5
$ws.Rows.Item(2).RowHeight = 19.5
$ws.Rows.Item(3).RowHeight = 19.5

$wb.Save()
